# "changes since last meeting" - add two new requirement rows to the
# Phase 1 Requirements sheet and move the active selection down to
# where the next entry would go.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = "Anjular JS, IONIX, JAVA, Spring Boot, Tomcat/Jetty, "
$ws.Range("B28").Value = "Sprint Management"

# Scroll the view down toward the new rows and leave the cursor where the
# user would type the next line.
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
[void]$ws.Range("B29").Select()
